# Updated cryptos list (Price / Volume(1h) columns) per upstream refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.212.53"
$ws.Range("E2").Value = "  +3.12%  "
$ws.Range("D3").Value = "2.995.82"
$ws.Range("E3").Value = "  +3.22%  "
$ws.Range("E4").Value = "  +0.35%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "562.25"
$ws.Range("E5").Value = "  +2.03%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "138.50"
$ws.Range("E6").Value = "  +12.35%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  +4.68%  "
$ws.Range("D9").Value = "2.988.23"
$ws.Range("E9").Value = "  +3.11%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.132"
$ws.Range("E10").Value = "  +7.26%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "4.94"
$ws.Range("E11").Value = "  +5.01%  "
$ws.Range("E12").Value = "  +4.06%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000230"
$ws.Range("E13").Value = "  +8.68%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "33.79"
$ws.Range("E14").Value = "  +3.88%  "
$ws.Range("E15").Value = "  +2.80%  "
$ws.Range("D16").Value = "3.488.57"
$ws.Range("E16").Value = "  +3.27%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.00"
$ws.Range("E17").Value = "  +7.13%  "
$ws.Range("D18").Value = "2.989.43"
$ws.Range("E18").Value = "  +3.24%  "
$ws.Range("D19").Value = "59.215.73"
$ws.Range("E19").Value = "  +3.21%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "424.85"
$ws.Range("E20").Value = "  +5.21%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.54"
$ws.Range("E21").Value = "  +5.03%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.711"
$ws.Range("E22").Value = "  +6.04%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.14"
$ws.Range("E23").Value = "  +4.50%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.43"
$ws.Range("E24").Value = "  +5.19%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "80.53"
$ws.Range("E25").Value = "  +4.45%  "
$ws.Range("E26").Value = "  +0.00%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("E27").Value = "  +0.54%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.13"
$ws.Range("E28").Value = "  +10.17%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.53"
$ws.Range("E29").Value = "  +3.34%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.81"
$ws.Range("E30").Value = "  +8.88%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "25.64"
$ws.Range("E31").Value = "  +3.89%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.14"
$ws.Range("E32").Value = "  +2.03%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0991"
$ws.Range("E33").Value = "  +0.22%  "
$ws.Range("E34").Value = "  +11.27%  "
$ws.Range("D35").Value = "0.0₃0770"
$ws.Range("E35").Value = "  +24.37%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.75"
$ws.Range("E36").Value = "  +6.09%  "
$ws.Range("E37").Value = "  +4.81%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "49.05"
$ws.Range("E38").Value = "  +2.32%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.66"
$ws.Range("E39").Value = "  +4.26%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.80"
$ws.Range("E40").Value = "  +17.06%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "402.59"
$ws.Range("E41").Value = "  +11.76%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0350"
$ws.Range("E42").Value = "  +2.69%  "
$ws.Range("D43").Value = "2.744.04"
$ws.Range("E43").Value = "  +4.69%  "
$ws.Range("E44").Value = "  +1.56%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.246"
$ws.Range("E45").Value = "  +7.49%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "125.24"
$ws.Range("E47").Value = "  +5.41%  "
$ws.Range("E49").Value = "  +2.47%  "
$ws.Range("E50").Value = "  +20.81%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "23.48"
$ws.Range("E51").Value = "  +2.86%  "
